$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diary already has an entry in row 14 (date / time start / time end /
# members present / discussion). Two further meetings are being logged in
# rows 16 and 18, following the same blank-row-separated layout and reusing
# the same cell formatting as the previous entry (row 14).

$ws.Range("A14:E14").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null

$ws.Range("A14:E14").Copy() | Out-Null
$ws.Range("A18:E18").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Meeting on 5 Oct 2023, 8pm - 10pm
$ws.Range("A16").Value = 45204
$ws.Range("B16").Value = 0.83333333333333337
$ws.Range("C16").Value = 0.91666666666666663
$ws.Range("D16").Value = "All"
$ws.Range("E16").Value = "Discussing temperature data and what to plot"

# Meeting on 6 Oct 2023, 11am - 4:30pm
$ws.Range("A18").Value = 45205
$ws.Range("B18").Value = 0.45833333333333331
$ws.Range("C18").Value = 0.6875
$ws.Range("D18").Value = "All "
$ws.Range("E18").Value = "Finalize the report and fixing aesthestics "
